$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "Dia Atual: 16/04/2023" -> "Dia Atual: 19/04/2023"
#         (only the run containing the lone "6" changes to "9"; the
#         surrounding runs " 1", "/04", "/2023" must stay untouched)
# ---------------------------------------------------------------------

# Anchor on the literal date text so we don't depend on fixed offsets.
$dateRng = $d.Content
$dateRng.Find.Execute("16/04/2023", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateStart = $dateRng.Start   # position of the "1" in "16"

# Segment ranges relative to $dateStart:
$rPrefix  = $d.Range($dateStart - 1, $dateStart + 1)   # " 1"
$rSix     = $d.Range($dateStart + 1, $dateStart + 2)   # "6"  (-> "9")
$rSlash04 = $d.Range($dateStart + 2, $dateStart + 5)   # "/04"
$rTail    = $d.Range($dateStart + 5, $dateStart + 10)  # "/2023"

# Temporarily give the neighbouring runs distinct underline formatting so
# that rewriting $rSix's text doesn't get coalesced into one big run with
# its (identically-formatted) neighbours.
$rPrefix.Font.Underline = 0
$rSlash04.Font.Underline = 3
$rTail.Font.Underline = 4

$rSix.Text = "9"

# Restore the neighbours' formatting (property-only changes do not trigger
# the run-coalescing pass that a text edit does).
$d.Range($dateStart - 1, $dateStart + 1).Font.Underline = 1
$d.Range($dateStart + 2, $dateStart + 5).Font.Underline = 1
$d.Range($dateStart + 5, $dateStart + 10).Font.Underline = 1

# ---------------------------------------------------------------------
# Edit 2: "Bloco Atual: Dark Hour" -> "Bloco Atual: Noite"
#         (leading space stays in its own run; "Dark Hour" is replaced by
#         a new run containing "Noite")
# ---------------------------------------------------------------------

$blockRng = $d.Content
$blockRng.Find.Execute("Dark Hour", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$blockRng.Text = ""
$blockRng.InsertAfter("Noite")
$blockRng.Font.Underline = 1
